$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting A:E to B:F
$ws.Columns("A:A").Insert()

# Copy the header formatting (bold, border, centered) from the shifted
# header in B1 onto the new A1 header cell, then set its text.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "ID"

# Populate the new ID column with the row identifiers
$ids = @{
    2 = "Hb 2"
    3 = "Hb 3"
    4 = "S 24"
    5 = "S 28"
    6 = "Hb 107"
    7 = "Hb 66"
    8 = "Hb 69"
    9 = "Hb 95"
    10 = "Hb 99"
    11 = "Hb 92"
    12 = "Hb 40"
    13 = "Hb 41"
    14 = "S 11"
    15 = "Hb 57"
    16 = "S 21"
    17 = "S 22"
    18 = "S 3"
    19 = "S 4"
    20 = "S 5"
    21 = "Hb 74"
    22 = "Hb 79"
    23 = "Hb 32"
    24 = "S 15"
    25 = "S 16"
}

foreach ($row in $ids.Keys) {
    $ws.Cells.Item($row, 1).Value = $ids[$row]
}
